# Applies the diff:
#  - Slide 2: merge title runs "What: The " + "Smart Player API" -> "What: The Smart Player API"
#             bump body placeholder font size 30 -> 32 (sz 3000 -> 3200) on the
#             visible text only (leave the two trailing empty paragraphs at sz 3000)
#  - Slide 3: bump URL placeholder font size 30 -> 32
#  - Slide 4: bump agenda placeholder font size 30 -> 32
#  - Slide 5: merge title runs "Smart Player API " + "Training & Pre-" -> "Smart Player API Training & Pre-"
#             bump body placeholder font size 30 -> 32, and merge runs
#             "basic HTML and JavaScript " + "experience" -> "basic HTML and JavaScript experience"

$p = $ppt.ActivePresentation

# Helper: the effective (visible) length of a TextRange, i.e. excluding
# trailing paragraph-mark-only (CR) characters from empty trailing paragraphs.
function Get-EffectiveLength($tr) {
    $effLen = $tr.Length
    while ($effLen -gt 0) {
        $ch = $tr.Characters($effLen, 1).Text
        if ($ch.Length -gt 0 -and [int][char]$ch[0] -eq 13) {
            $effLen = $effLen - 1
        } else {
            break
        }
    }
    return $effLen
}

# ---------------------------------------------------------------------------
# Slide 2
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Title: "What: The " + "Smart Player API" -> single run
# (round-trip through a placeholder value first so the COM layer actually
# rewrites/merges the runs instead of treating an identical string as a no-op)
$title2 = $s2.Shapes.Item(2)
$title2.TextFrame.TextRange.Text = "__tmp__"
$title2.TextFrame.TextRange.Text = "What: The Smart Player API"

# Body: bump font size of the visible text (sz 3000 -> 3200)
$body2 = $s2.Shapes.Item(3)
$tr2 = $body2.TextFrame.TextRange
$eff2 = Get-EffectiveLength($tr2)
$tr2.Characters(1, $eff2).Font.Size = 32

# ---------------------------------------------------------------------------
# Slide 3
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2)
$tr3 = $body3.TextFrame.TextRange
$eff3 = Get-EffectiveLength($tr3)
$tr3.Characters(1, $eff3).Font.Size = 32

# ---------------------------------------------------------------------------
# Slide 4
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$body4 = $s4.Shapes.Item(2)
$tr4 = $body4.TextFrame.TextRange
$eff4 = Get-EffectiveLength($tr4)
$tr4.Characters(1, $eff4).Font.Size = 32

# ---------------------------------------------------------------------------
# Slide 5
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)

# Title: "Smart Player API " + "Training & Pre-" -> single run (keep "reqs" run separate)
$title5 = $s5.Shapes.Item(2)
$trTitle5 = $title5.TextFrame.TextRange
$trTitle5.Characters(1, 32).Text = "Smart Player API Training & Pre-"

# Body: bump font size, and merge "basic HTML and JavaScript " + "experience"
$body5 = $s5.Shapes.Item(3)
$tr5 = $body5.TextFrame.TextRange
$eff5 = Get-EffectiveLength($tr5)
$tr5.Characters(1, $eff5).Font.Size = 32
$tr5.Characters(167, 36).Text = "basic HTML and JavaScript experience"
